$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Inscritos)
$ws.Range("E8").Value = 36

# Row 15 (Inscritos, Pagos, Inscrições homologadas)
$ws.Range("E15").Value = 82
$ws.Range("F15").Value = 39
$ws.Range("H15").Value = 39

# Row 16 (Inscritos)
$ws.Range("E16").Value = 280

# Row 18 (Inscritos, Pagos, Inscrições homologadas)
$ws.Range("E18").Value = 83
$ws.Range("F18").Value = 27
$ws.Range("H18").Value = 27
